$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "17/04/2020" case count (row 36): 797 -> 920
$ws.Range("B36").Value = 920

# Insert a new row for "18/04/2020" right after "18/03/2020" (row 37),
# pushing "19/03/2020" (previously row 38) and everything after it down by one.
$ws.Range("A38").EntireRow.Insert()

# Populate the newly inserted row 38 with the new data point.
$ws.Range("A38").Value = "18/04/2020"
$ws.Range("B38").Value = 324
